# Insert a new weekly price-report row for "Apio" (Vega Monumental Concepción)
# at sheet row 112, pushing the existing rows 112:145 down to 113:146.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("112:112").Insert()

$ws.Range("A112").Value = 11
$ws.Range("B112").Value = "Vega Monumental Concepción"
$ws.Range("C112").Value = "Bíobío"
$ws.Range("D112").Value = 44463
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 100112017
$ws.Range("G112").Value = "Apio"
$ws.Range("H112").Value = "Americana (o)"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 100
$ws.Range("K112").Value = 8000
$ws.Range("L112").Value = 8500
$ws.Range("M112").Value = 8250
$ws.Range("N112").Value = "$/docena de matas"
$ws.Range("O112").Value = "Región de Coquimbo"
$ws.Range("P112").Value = 1375
$ws.Range("Q112").Value = 6
$ws.Range("R112").Value = "Hortaliza"
